# Discharge_Aug07.xlsx -- "lots of discharge data"
#
# Adds a third depth/velocity/discharge block (rows 29-40) to the "stn3"
# sheet, mirroring the two existing blocks (rows 3-12 and rows 15-26) but
# with column C ("D", i.e. the gauge-height/segment reading) recomputed in
# centimeters from the previous block (C17:C26 * 2.54), and column B
# ("V", velocity) pasted in as plain values copied from the B17:B26 block.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "stn1"
$ws3 = $wb.Worksheets.Item(2)   # "stn3"

# ---------------------------------------------------------------------
# New data block on stn3 (rows 29-40)
# ---------------------------------------------------------------------

# Row 29: bold "new depth" section header (reuses the same shared string /
# style already used for the "new velocity" header above it on row 15).
$ws3.Range("A29").Value = "new depth"
$ws3.Range("A29").Font.Bold = $true

# Row 30: column headers, same as row 2 / row 16 above.
$ws3.Range("A30").Value = "X"
$ws3.Range("B30").Value = "V"
$ws3.Range("C30").Value = "D"
$ws3.Range("D30").Value = "segment"
$ws3.Range("E30").Value = "Q"
$ws3.Range("F30").Value = "Qtotal"

# Row 31-40, column A: depth values (same series as A17:A26).
$aVals = @(0.2, 0.25, 0.3, 0.35, 0.4, 0.45, 0.5, 0.55, 0.6, 0.65)
for ($i = 0; $i -lt 10; $i++) {
    $r = 31 + $i
    $ws3.Range("A$r").Value = $aVals[$i]
}

# Row 31-40, column B: velocity -- pasted as plain VALUES copied from the
# B17:B26 block above (not a formula).
for ($i = 0; $i -lt 10; $i++) {
    $r = 31 + $i
    $srcRow = 17 + $i
    $ws3.Range("B$r").Value = $ws3.Range("B$srcRow").Value2
}

# Row 31-40, column C: "D" reading converted to centimeters from the
# corresponding row in the block above (C17:C26 * 2.54).
$ws3.Range("C31").Formula = "=C17*2.54"
$ws3.Range("C32:C40").Formula = "=C18*2.54"

# Row 31, column D: same as D17/D3 -- just echoes A31.
$ws3.Range("D31").Formula = "=A31"
# Row 32, column D: midpoint formula (first of its run, like D18/D4).
$ws3.Range("D32").Formula = "=(A32+(A33-A32)/2)"
# Rows 33-40, column D: shared midpoint formula.
$ws3.Range("D33:D40").Formula = "=(A33+(A34-A33)/2)"

# Row 32, column E: segment discharge (first of its run, like E18/E4).
$ws3.Range("E32").Formula = "=(D32-D31)*(B32)*C32"
# Rows 33-40, column E: shared segment discharge formula.
$ws3.Range("E33:E40").Formula = "=(D33-D32)*(B33)*C33"

# Row 31, column F: total discharge for the block.
$ws3.Range("F31").Formula = "=SUM(E31:E40)"

# ---------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------

# stn1 becomes the non-selected tab, scrolled/selected at C42.
[void]$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("C42").Select()

# stn3 ends up the active/selected tab, scrolled to row 3, selection on
# the new header cell A29.
[void]$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws3.Range("A29").Select()
